$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 37039140
$ws.Range("I40").Value = 1949.75
$ws.Range("J40").Value = 43480390
$ws.Range("K40").Value = 1949.75
$ws.Range("L40").Value = 43480390
$ws.Range("M40").Value = -1774.75
$ws.Range("N40").Value = -43480740
$ws.Range("H41").Value = 720.1
$ws.Range("I41").Value = 205
$ws.Range("J41").Value = 1063.5
$ws.Range("K41").Value = 205
$ws.Range("L41").Value = 1063.5
$ws.Range("M41").Value = 235
$ws.Range("N41").Value = -1943.5
$ws.Range("H112").Value = 1490.641
$ws.Range("I112").Value = 993.625
$ws.Range("J112").Value = 1618.9032
$ws.Range("K112").Value = 2980.875
$ws.Range("L112").Value = 4856.7096
$ws.Range("M112").Value = -1872.875
$ws.Range("N112").Value = -7072.7096
$ws.Range("H113").Value = 2524.147
$ws.Range("I113").Value = 2759.7058
$ws.Range("J113").Value = 2288.5881
$ws.Range("K113").Value = 2759.7058
$ws.Range("L113").Value = 2288.5881
$ws.Range("M113").Value = 494.2941999999998
$ws.Range("N113").Value = -8796.588100000001
$ws.Range("H116").Value = 4526
$ws.Range("I116").Value = 4628.4
$ws.Range("J116").Value = 4270
$ws.Range("K116").Value = 4628.4
$ws.Range("L116").Value = 4270
$ws.Range("M116").Value = -1186.4
$ws.Range("N116").Value = -11154
$ws.Range("H132").Value = 4894.44
$ws.Range("I132").Value = 4970.9546
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 14912.8638
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -12382.8638
$ws.Range("N132").Value = -18060.0005
$ws.Range("H138").Value = 1451.14
$ws.Range("I138").Value = 639.9219000000001
$ws.Range("J138").Value = 2893.3057
$ws.Range("K138").Value = 1919.7657
$ws.Range("L138").Value = 8679.917099999999
$ws.Range("M138").Value = 3220.2343
$ws.Range("N138").Value = -18959.9171
$ws.Range("H141").Value = 2612.1482
$ws.Range("I141").Value = 670.2895
$ws.Range("J141").Value = 7224.0625
$ws.Range("K141").Value = 2010.8685
$ws.Range("L141").Value = 21672.1875
$ws.Range("M141").Value = 3169.1315
$ws.Range("N141").Value = -32032.1875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1051.375
$ws.Range("I61").Value = 755.53845
$ws.Range("J61").Value = 2333.3333
$ws.Range("K61").Value = 755.53845
$ws.Range("L61").Value = 2333.3333
$ws.Range("M61").Value = -543.53845
$ws.Range("N61").Value = -2757.3333
$ws.Range("H74").Value = 738
$ws.Range("I74").Value = 686.4706
$ws.Range("J74").Value = 1000.8
$ws.Range("K74").Value = 686.4706
$ws.Range("L74").Value = 1000.8
$ws.Range("M74").Value = 187.5294
$ws.Range("N74").Value = -2748.8
$ws.Range("H77").Value = 738
$ws.Range("I77").Value = 686.4706
$ws.Range("J77").Value = 1000.8
$ws.Range("K77").Value = 3432.353
$ws.Range("L77").Value = 5004
$ws.Range("M77").Value = 935.6469999999999
$ws.Range("N77").Value = -13740
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H132").Value = 1589
$ws.Range("I132").Value = 1157.0714
$ws.Range("J132").Value = 3604.6667
$ws.Range("K132").Value = 3471.2142
$ws.Range("L132").Value = 10814.0001
$ws.Range("M132").Value = -941.2142000000003
$ws.Range("N132").Value = -15874.0001
$ws.Range("H136").Value = 1051.375
$ws.Range("I136").Value = 755.53845
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 2266.61535
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = 283.38465
$ws.Range("N136").Value = -12099.9999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 45456212
$ws.Range("I99").Value = 62501628
$ws.Range("J99").Value = 1770.3334
$ws.Range("K99").Value = 62501628
$ws.Range("L99").Value = 1770.3334
$ws.Range("M99").Value = -62500130
$ws.Range("N99").Value = -4766.3334
$ws.Range("H134").Value = 22081.256
$ws.Range("I134").Value = 2209.359
$ws.Range("J134").Value = 86664.914
$ws.Range("K134").Value = 6628.076999999999
$ws.Range("L134").Value = 259994.742
$ws.Range("M134").Value = -4093.076999999999
$ws.Range("N134").Value = -265064.742

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2125.149
$ws.Range("I31").Value = 2155.1555
$ws.Range("J31").Value = 1450
$ws.Range("K31").Value = 2155.1555
$ws.Range("L31").Value = 1450
$ws.Range("M31").Value = -1860.1555
$ws.Range("N31").Value = -2040
$ws.Range("H34").Value = 2125.149
$ws.Range("I34").Value = 2155.1555
$ws.Range("J34").Value = 1450
$ws.Range("K34").Value = 2155.1555
$ws.Range("L34").Value = 1450
$ws.Range("M34").Value = -1953.1555
$ws.Range("N34").Value = -1854
$ws.Range("H99").Value = 2430.4473
$ws.Range("I99").Value = 2089.1304
$ws.Range("J99").Value = 2953.8
$ws.Range("K99").Value = 2089.1304
$ws.Range("L99").Value = 2953.8
$ws.Range("M99").Value = -591.1304
$ws.Range("N99").Value = -5949.8
$ws.Range("H126").Value = 2430.4473
$ws.Range("I126").Value = 2089.1304
$ws.Range("J126").Value = 2953.8
$ws.Range("K126").Value = 6267.3912
$ws.Range("L126").Value = 8861.400000000001
$ws.Range("M126").Value = -3797.3912
$ws.Range("N126").Value = -13801.4
$ws.Range("H132").Value = 3009
$ws.Range("I132").Value = 2340
$ws.Range("J132").Value = 3566.5
$ws.Range("K132").Value = 7020
$ws.Range("L132").Value = 10699.5
$ws.Range("M132").Value = -4490
$ws.Range("N132").Value = -15759.5
$ws.Range("H134").Value = 1646.4615
$ws.Range("I134").Value = 1044.8889
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3134.6667
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -599.6666999999998
$ws.Range("N134").Value = -14070

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1851.138
$ws.Range("I4").Value = 65.5
$ws.Range("J4").Value = 2316.9565
$ws.Range("K4").Value = 196.5
$ws.Range("L4").Value = 6950.869499999999
$ws.Range("M4").Value = -84.5
$ws.Range("N4").Value = -7174.869499999999
$ws.Range("H14").Value = 129.4
$ws.Range("I14").Value = 129.4
$ws.Range("K14").Value = 388.2
$ws.Range("M14").Value = -215.2
$ws.Range("H92").Value = 1086
$ws.Range("I92").Value = 267
$ws.Range("J92").Value = 6000
$ws.Range("K92").Value = 801
$ws.Range("L92").Value = 18000
$ws.Range("M92").Value = 447
$ws.Range("N92").Value = -20496
$ws.Range("H107").Value = 556141.0600000001
$ws.Range("I107").Value = 1112.6
$ws.Range("J107").Value = 864490.25
$ws.Range("K107").Value = 3337.8
$ws.Range("L107").Value = 2593470.75
$ws.Range("M107").Value = -1417.8
$ws.Range("N107").Value = -2597310.75
$ws.Range("H114").Value = 1179.5
$ws.Range("I114").Value = 621.125
$ws.Range("J114").Value = 1737.875
$ws.Range("K114").Value = 1863.375
$ws.Range("L114").Value = 5213.625
$ws.Range("M114").Value = 1390.625
$ws.Range("N114").Value = -11721.625
$ws.Range("H116").Value = 127103.5
$ws.Range("I116").Value = 1138.1666
$ws.Range("J116").Value = 504999.5
$ws.Range("K116").Value = 3414.4998
$ws.Range("L116").Value = 1514998.5
$ws.Range("M116").Value = 27.50019999999995
$ws.Range("N116").Value = -1521882.5
$ws.Range("H131").Value = 826.4
$ws.Range("J131").Value = 863.337
$ws.Range("L131").Value = 2590.011
$ws.Range("N131").Value = -12670.011

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2585
$ws.Range("I132").Value = 2374.7058
$ws.Range("K132").Value = 7124.117400000001
$ws.Range("M132").Value = -4594.117400000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1168.0667
$ws.Range("I46").Value = 1135.0834
$ws.Range("K46").Value = 1135.0834
$ws.Range("M46").Value = -947.0834
$ws.Range("H132").Value = 2040.4
$ws.Range("I132").Value = 1648.6538
$ws.Range("J132").Value = 2926.087
$ws.Range("K132").Value = 4945.9614
$ws.Range("L132").Value = 8778.261
$ws.Range("M132").Value = -2415.9614
$ws.Range("N132").Value = -13838.261
$ws.Range("H136").Value = 3289.9167
$ws.Range("I136").Value = 1491.6471
$ws.Range("J136").Value = 7657.143
$ws.Range("K136").Value = 4474.9413
$ws.Range("L136").Value = 22971.429
$ws.Range("M136").Value = -1924.9413
$ws.Range("N136").Value = -28071.429

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8473
$ws.Range("J74").Value = 8473
$ws.Range("L74").Value = 8473
$ws.Range("N74").Value = -10345
$ws.Range("H77").Value = 8473
$ws.Range("J77").Value = 8473
$ws.Range("L77").Value = 25419
$ws.Range("N77").Value = -34779
$ws.Range("H113").Value = 748.2222
$ws.Range("I113").Value = 619.1429000000001
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1857.4287
$ws.Range("L113").Value = 3600
$ws.Range("M113").Value = 312.5712999999998
$ws.Range("N113").Value = -7940
$ws.Range("H136").Value = 725.64703
$ws.Range("I136").Value = 389.06668
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 1167.20004
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = 1382.79996
$ws.Range("N136").Value = -14850
